# toy_model4.xlsx test-file update
#
# Restores two metabolites ("m2" and "m8") that were missing from the
# metsData sheet's bound table, and refreshes the workbook/sheet view
# state (active sheet + selections) that LibreOffice re-wrote when the
# file was last saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. metsData: insert the two missing metabolite rows so the A column
#    again lists m1..m25 in order (m2 was missing after m1/m3, m8 was
#    missing after m7). Every other column for these rows follows the
#    existing 0.99 / 1 / 1.01 pattern used throughout the sheet.
# ---------------------------------------------------------------------
$metsData = $wb.Worksheets.Item("metsData")

$metsData.Rows(4).Insert()
$metsData.Range("A4").Value = "m2"
$metsData.Range("B4").Value = 0.99
$metsData.Range("C4").Value = 1
$metsData.Range("D4").Value = 1.01

$metsData.Rows(9).Insert()
$metsData.Range("A9").Value = "m8"
$metsData.Range("B9").Value = 0.99
$metsData.Range("C9").Value = 1
$metsData.Range("D9").Value = 1.01

# ---------------------------------------------------------------------
# 2. Column width touch-ups left behind by the resave (best effort;
#    Excel snaps these to its internal character-width granularity).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("general").Columns("A").ColumnWidth = 73.1619433198381
$wb.Worksheets.Item("general").Columns("B").ColumnWidth = 30.7449392712551

$wb.Worksheets.Item("measRates").Columns("A").ColumnWidth = 20.1376518218623
$wb.Worksheets.Item("measRates").Columns("C").ColumnWidth = 13.3886639676113

$wb.Worksheets.Item("kinetics1").Columns("A").ColumnWidth = 15.7449392712551
$wb.Worksheets.Item("kinetics1").Columns("B").ColumnWidth = 34.4939271255061
$wb.Worksheets.Item("kinetics1").Columns("C").ColumnWidth = 12.6396761133603

$wb.Worksheets.Item("stoic").Columns("A").ColumnWidth = 12.748987854251

$wb.Worksheets.Item("thermoRxns").Columns("B").ColumnWidth = 16.2834008097166

# ---------------------------------------------------------------------
# 3. Re-select each sheet's used columns (A:D), which is what every
#    sheet's saved selection became.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("general").Range("A:D").Select()
$wb.Worksheets.Item("stoic").Range("A:D").Select()
$wb.Worksheets.Item("mets").Range("A:D").Select()
$wb.Worksheets.Item("rxns").Range("A:D").Select()
$wb.Worksheets.Item("splitRatios").Range("A:D").Select()
$wb.Worksheets.Item("poolConst").Range("A:D").Select()
$wb.Worksheets.Item("thermo_ineq_constraints").Range("A:D").Select()
$wb.Worksheets.Item("thermoRxns").Range("A:D").Select()
$wb.Worksheets.Item("thermoMets").Range("A:D").Select()
$wb.Worksheets.Item("measRates").Range("A:D").Select()
$wb.Worksheets.Item("protData").Range("A:D").Select()
$wb.Worksheets.Item("kinetics1").Range("A:D").Select()

# metsData's selection becomes A1 / A:D.
$metsData.Range("A:D").Select()

# ---------------------------------------------------------------------
# 4. metsData is the sheet left active/selected when the workbook was
#    saved (activeTab moves from "general" (0) to "metsData" (11)).
# ---------------------------------------------------------------------
$metsData.Activate()
